$d = $word.ActiveDocument

# The paragraph originally read (across several runs):
#   " ... 300 frames for all replacement policies, supporting the idea that
#   rand in general isn't as bad when working with programs with smaller
#   memory needs."
# The commit trims it down to:
#   " ... 300 frames for all replacement policies."
# i.e. the clause starting at the comma after "policies" through
# "...smaller memory needs" is removed, leaving the trailing period intact.

$old = " for all replacement policies, supporting the idea that rand in general isn" + [char]0x2019 + "t as bad when working with programs with smaller memory needs"
$new = " for all replacement policies"

$range = $d.Content
$range.Find.ClearFormatting()
$found = $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)

if (-not $found) {
    throw "edit.ps1: target sentence fragment not found; document may not match expected starting state"
}
